# Update the transit-depth data table: replace the wavelength (B) and
# D_lambda (C) series with the new, finer-resolution values used for the
# paper figures, extending the table from 49 to 58 data rows
# (A1:C50 -> A1:C59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement dataset (58 rows): column A is the 0-based sample index,
# column B is wavelength [nm], column C is D_lambda [ppm].
$A = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57)
$B = @(300.5,314.5,328.4,342,356,369.6,383.5,397.5,411.1,425,439,452.6,466.6,480.5,494.1,508.3,522,535.5999999999999,550,563.5,588.4000000000001,614,639,664.6,689.5,715,740,765,790.5,816,841.3000000000001,866,891.5999999999999,917,942.4,967.3000000000001,992.7,1018,1055,1083,1112,1143,1174,1205,1237,1271,1306,1341,1377,1414,1452,1492,1532,1573,1616,1659,1704,1751)
$C = @(13829.57077530367,13798.34573850303,13784.82275098525,13757.7111932744,13733.50841330345,13710.96301541541,13729.22856815839,13675.52039284026,13685.15160235384,13651.03537687185,13632.27689882129,13615.56496736183,13592.55321879471,13576.07599865851,13545.1209805646,13525.25262210325,13494.53871429185,13479.39767999029,13464.08422607948,13443.66744424574,13412.99572174015,13378.40109863953,13348.14969893239,13317.25640898218,13290.83580293278,13262.42358578056,13236.51124327019,13212.54115925485,13191.19526095858,13169.59761725967,13147.19450890778,13129.47389008012,13107.60827624802,13093.82921692181,13077.81352005033,13059.79395945012,13051.16955260255,13036.18844731491,13014.84125138863,12997.62738210208,12989.09665289716,12973.88741734184,12960.65204192398,12950.86462541251,12938.62893550024,12925.86418645381,12915.17140127718,12900.73781980283,12888.04447691916,12872.7927092902,12861.4395289921,12846.26581892168,12830.77116018516,12814.85227933144,12798.58931511302,12793.46200801357,12777.35696655291,12763.48093909863)

# Reference the existing style applied to column A (bold/bordered header
# style, already used on A2 in the original sheet) so that newly created
# rows 51-59 keep the same look as the pre-existing ones.
$aStyle = $ws.Range("A2").Style

for ($i = 0; $i -lt $A.Length; $i++) {
    $row = $i + 2
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $A[$i]
    $cellA.Style = $aStyle
    $ws.Cells.Item($row, 2).Value = $B[$i]
    $ws.Cells.Item($row, 3).Value = $C[$i]
}
